$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Divya Bhagas: E2 1->3, F2 100->33.3, G2 "91-100%"->"26-50%"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 33.3
$ws.Range("G2").Value = "26-50%"

# Row 3 - Sai Kenekar: E3 0->2
$ws.Range("E3").Value = 2

# Row 4 - Anushka Mote: E4 1->3, F4 100->33.3, G4 "91-100%"->"26-50%"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 33.3
$ws.Range("G4").Value = "26-50%"

# Row 5 - Abhishek Pathak: E5 0->2
$ws.Range("E5").Value = 2

# Row 6 - Vaishnavi Pawar: E6 0->2
$ws.Range("E6").Value = 2

# Row 7 - Sagar Pawar: E7 2->4, F7 100->50, G7 "91-100%"->"26-50%"
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 50
$ws.Range("G7").Value = "26-50%"

# Row 8 - Shubham Phad: E8 1->3, F8 100->33.3, G8 "91-100%"->"26-50%"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 33.3
$ws.Range("G8").Value = "26-50%"

# Row 9 - Shubham Pitekar: D9 4->6, E9 4->6 (F9 stays 100, G9 stays "91-100%")
$ws.Range("D9").Value = 6
$ws.Range("E9").Value = 6

# Row 10 - Damini Solunke: E10 1->3, F10 100->33.3, G10 "91-100%"->"26-50%"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 33.3
$ws.Range("G10").Value = "26-50%"
